# Refresh NATMI ligand-receptor (Ntn1 -> Adora2b) edge stats with the
# re-run TPM-based scRNA-seq expression numbers (ligand-/receptor-expressing
# cell counts, detection rates, average/total expression, derived
# specificities, and edge weights) for every Sending x Target cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.532141
$ws.Range("H2").Value = 4.596423
$ws.Range("I2").Value = 0.08900664250669833
$ws.Range("J2").Value = 0.08900664250669831
$ws.Range("O2").Value = 0.4853706552224367
$ws.Range("P2").Value = 0.4853706552224367
$ws.Range("Q2").Value = 5.958141403001666
$ws.Range("R2").Value = 53.62327262701499
$ws.Range("S2").Value = 0.04320121239262536
$ws.Range("T2").Value = 0.04320121239262534

# Row 3: ECs -> MuSCs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.532141
$ws.Range("H3").Value = 4.596423
$ws.Range("I3").Value = 0.08900664250669833
$ws.Range("J3").Value = 0.08900664250669831
$ws.Range("M3").Value = 2.316983333333333
$ws.Range("N3").Value = 6.950950000000001
$ws.Range("O3").Value = 0.289190721133932
$ws.Range("P3").Value = 0.289190721133932
$ws.Range("Q3").Value = 3.549945161316667
$ws.Range("R3").Value = 31.94950645185
$ws.Range("S3").Value = 0.02573989513222218
$ws.Range("T3").Value = 0.02573989513222217

# Row 4: ECs -> Resolving-Mac
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.532141
$ws.Range("H4").Value = 4.596423
$ws.Range("I4").Value = 0.08900664250669833
$ws.Range("J4").Value = 0.08900664250669831
$ws.Range("M4").Value = 1.806204333333334
$ws.Range("N4").Value = 5.418613000000001
$ws.Range("O4").Value = 0.2254386236436313
$ws.Range("P4").Value = 0.2254386236436313
$ws.Range("Q4").Value = 2.767359713477667
$ws.Range("R4").Value = 24.906237421299
$ws.Range("S4").Value = 0.0200655349818508
$ws.Range("T4").Value = 0.02006553498185079

# Row 5: FAPs -> FAPs
$ws.Range("I5").Value = 0.6169137955113024
$ws.Range("J5").Value = 0.6169137955113023
$ws.Range("O5").Value = 0.4853706552224367
$ws.Range("P5").Value = 0.4853706552224367
$ws.Range("S5").Value = 0.2994318531430812
$ws.Range("T5").Value = 0.2994318531430811

# Row 6: FAPs -> MuSCs
$ws.Range("I6").Value = 0.6169137955113024
$ws.Range("J6").Value = 0.6169137955113023
$ws.Range("M6").Value = 2.316983333333333
$ws.Range("N6").Value = 6.950950000000001
$ws.Range("O6").Value = 0.289190721133932
$ws.Range("P6").Value = 0.289190721133932
$ws.Range("Q6").Value = 24.60501914966667
$ws.Range("R6").Value = 221.445172347
$ws.Range("S6").Value = 0.1784057454013846
$ws.Range("T6").Value = 0.1784057454013846

# Row 7: FAPs -> Resolving-Mac
$ws.Range("I7").Value = 0.6169137955113024
$ws.Range("J7").Value = 0.6169137955113023
$ws.Range("M7").Value = 1.806204333333334
$ws.Range("N7").Value = 5.418613000000001
$ws.Range("O7").Value = 0.2254386236436313
$ws.Range("P7").Value = 0.2254386236436313
$ws.Range("Q7").Value = 19.18084242148667
$ws.Range("R7").Value = 172.62758179338
$ws.Range("S7").Value = 0.1390761969668366
$ws.Range("T7").Value = 0.1390761969668366

# Row 8: MuSCs -> FAPs
$ws.Range("G8").Value = 4.902263666666666
$ws.Range("H8").Value = 14.706791
$ws.Range("I8").Value = 0.2847871244569372
$ws.Range("J8").Value = 0.2847871244569371
$ws.Range("O8").Value = 0.4853706552224367
$ws.Range("P8").Value = 0.4853706552224367
$ws.Range("Q8").Value = 19.06376770858389
$ws.Range("R8").Value = 171.573909377255
$ws.Range("S8").Value = 0.1382273131965772
$ws.Range("T8").Value = 0.1382273131965772

# Row 9: MuSCs -> MuSCs
$ws.Range("G9").Value = 4.902263666666666
$ws.Range("H9").Value = 14.706791
$ws.Range("I9").Value = 0.2847871244569372
$ws.Range("J9").Value = 0.2847871244569371
$ws.Range("M9").Value = 2.316983333333333
$ws.Range("N9").Value = 6.950950000000001
$ws.Range("O9").Value = 0.289190721133932
$ws.Range("P9").Value = 0.289190721133932
$ws.Range("Q9").Value = 11.35846321127222
$ws.Range("R9").Value = 102.22616890145
$ws.Range("S9").Value = 0.0823577938913605
$ws.Range("T9").Value = 0.08235779389136048

# Row 10: MuSCs -> Resolving-Mac
$ws.Range("G10").Value = 4.902263666666666
$ws.Range("H10").Value = 14.706791
$ws.Range("I10").Value = 0.2847871244569372
$ws.Range("J10").Value = 0.2847871244569371
$ws.Range("M10").Value = 1.806204333333334
$ws.Range("N10").Value = 5.418613000000001
$ws.Range("O10").Value = 0.2254386236436313
$ws.Range("P10").Value = 0.2254386236436313
$ws.Range("Q10").Value = 8.854489877875888
$ws.Range("R10").Value = 79.690408900883
$ws.Range("S10").Value = 0.06420201736899944
$ws.Range("T10").Value = 0.06420201736899941

# Row 11: Resolving-Mac -> FAPs
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.159958
$ws.Range("H11").Value = 0.479874
$ws.Range("I11").Value = 0.009292437525062282
$ws.Range("J11").Value = 0.009292437525062281
$ws.Range("O11").Value = 0.4853706552224367
$ws.Range("P11").Value = 0.4853706552224367
$ws.Range("Q11").Value = 0.6220396050633332
$ws.Range("R11").Value = 5.598356445569999
$ws.Range("S11").Value = 0.004510276490153039
$ws.Range("T11").Value = 0.004510276490153038

# Row 12: Resolving-Mac -> MuSCs
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.159958
$ws.Range("H12").Value = 0.479874
$ws.Range("I12").Value = 0.009292437525062282
$ws.Range("J12").Value = 0.009292437525062281
$ws.Range("M12").Value = 2.316983333333333
$ws.Range("N12").Value = 6.950950000000001
$ws.Range("O12").Value = 0.289190721133932
$ws.Range("P12").Value = 0.289190721133932
$ws.Range("Q12").Value = 0.3706200200333333
$ws.Range("R12").Value = 3.3355801803
$ws.Range("S12").Value = 0.002687286708964772
$ws.Range("T12").Value = 0.002687286708964772

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.159958
$ws.Range("H13").Value = 0.479874
$ws.Range("I13").Value = 0.009292437525062282
$ws.Range("J13").Value = 0.009292437525062281
$ws.Range("M13").Value = 1.806204333333334
$ws.Range("N13").Value = 5.418613000000001
$ws.Range("O13").Value = 0.2254386236436313
$ws.Range("P13").Value = 0.2254386236436313
$ws.Range("Q13").Value = 0.2889168327513333
$ws.Range("R13").Value = 2.600251494762
$ws.Range("S13").Value = 0.002094874325944472
$ws.Range("T13").Value = 0.002094874325944472
